$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = -22.14670000000001
$ws.Range("A21").Value = -20.06619999999998
$ws.Range("A23").Value = -20.55039999999997
$ws.Range("B24").Value = 5.439000000000004
$ws.Range("A25").Value = -21.62499999999999
$ws.Range("B28").Value = 6.014500000000002
$ws.Range("B36").Value = 9.377100000000008
$ws.Range("B45").Value = 5.4613
$ws.Range("B48").Value = 7.039800000000003
$ws.Range("B49").Value = 6.093999999999995
$ws.Range("B52").Value = 5.455399999999997
$ws.Range("A53").Value = -21.77409999999999
$ws.Range("B53").Value = 5.689799999999998
$ws.Range("B54").Value = 4.8755
$ws.Range("A57").Value = -22.41050000000002
$ws.Range("A59").Value = -21.85449999999999
$ws.Range("A69").Value = -21.63799999999999
$ws.Range("B70").Value = 7.2288
$ws.Range("A79").Value = -20.03800000000002
$ws.Range("A83").Value = -21.6604
$ws.Range("B86").Value = 5.016899999999998
$ws.Range("B87").Value = 5.425999999999997
$ws.Range("A93").Value = -21.18840000000002
$ws.Range("B101").Value = 5.165999999999999
